$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.528.90"
$ws.Range("E2").Value = "  +0.88%  "

# Row 3
$ws.Range("D3").Value = "2.481.35"
$ws.Range("E3").Value = "  +0.80%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.88"
$ws.Range("E5").Value = "  +0.70%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.40"
$ws.Range("E6").Value = "  -0.38%  "

# Row 7
$ws.Range("E7").Value = "  -1.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.22%  "

# Row 9
$ws.Range("E9").Value = "  +2.32%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.82"
$ws.Range("E10").Value = "  -0.97%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("E11").Value = "  +1.32%  "

# Row 12
$ws.Range("E12").Value = "  +2.35%  "

# Row 13
$ws.Range("D13").Value = "2.862.11"
$ws.Range("E13").Value = "  +0.75%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.86"
$ws.Range("E14").Value = "  -1.31%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.23"
$ws.Range("E15").Value = "  +9.72%  "

# Row 16
$ws.Range("D16").Value = "2.499.09"
$ws.Range("E16").Value = "  +1.62%  "

# Row 17
$ws.Range("E17").Value = "  -1.87%  "

# Row 18
$ws.Range("D18").Value = "41.534.85"
$ws.Range("E18").Value = "  +0.99%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.42"
$ws.Range("E19").Value = "  +2.52%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0940"
$ws.Range("E20").Value = "  +2.36%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.66"
$ws.Range("E21").Value = "  +4.84%  "

# Row 22
$ws.Range("E22").Value = "  +1.64%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.10"
$ws.Range("E23").Value = "  +0.87%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.70"
$ws.Range("E24").Value = "  -1.43%  "

# Row 25
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.96"
$ws.Range("E27").Value = "  +4.27%  "

# Row 28
$ws.Range("E28").Value = "  +0.94%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.66"
$ws.Range("E29").Value = "  +0.95%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.10"
$ws.Range("E30").Value = "  -0.28%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.42"
$ws.Range("E31").Value = "  +3.93%  "

# Row 32
$ws.Range("E32").Value = "  +0.02%  "

# Row 33
$ws.Range("E33").Value = "  +0.62%  "

# Row 34
$ws.Range("E34").Value = "  +2.34%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.35"
$ws.Range("E35").Value = "  +2.42%  "

# Row 36
$ws.Range("E36").Value = "  -8.29%  "

# Row 37
$ws.Range("E37").Value = "  +4.59%  "

# Row 38
$ws.Range("E38").Value = "  -2.98%  "

# Row 39
$ws.Range("E39").Value = "  -2.29%  "

# Row 40
$ws.Range("E40").Value = "  +0.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.11"

# Row 42
$ws.Range("E42").Value = "  -0.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.66"
$ws.Range("E43").Value = "  -1.85%  "

# Row 44
$ws.Range("D44").Value = "1.974.41"
$ws.Range("E44").Value = "  +0.10%  "

# Row 45
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("E46").Value = "  -1.78%  "

# Row 47
$ws.Range("E47").Value = "  +3.50%  "

# Row 48
$ws.Range("D48").Value = "2.719.89"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.76"
$ws.Range("E49").Value = "  +1.06%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.19"
$ws.Range("E50").Value = "  -1.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.64"
$ws.Range("E51").Value = "  -1.44%  "
